# Daily attendance processing - 2025-12-01 16:34:15
#
# The "Recorded By" column (G) lists the users who recorded/edited each
# attendance session, as a comma-separated string (e.g. "System, someone@example.com").
# This pass normalizes the ordering of that list for every data row by moving
# the first-listed recorder to the end of the list (a left rotation by one
# element). Cells that only contain a single recorder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $rest = $parts[1..($parts.Count - 1)]
        $rotated = $rest + $parts[0]
        $newVal = $rotated -join ", "
        $cell.Value2 = $newVal
    }
}
